$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; temporarily unprotect so the cells can be updated,
# then restore protection once the edits are made.
$ws.Unprotect()

# Update the confidential disclaimer date in A13: 2021-05-05 -> 2021-05-06
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for rows 2-10
$ws.Range("D2").Value = 0.0958095402228978
$ws.Range("E2").Value = -0.0274089414057741

$ws.Range("D3").Value = 0.1069224073228835
$ws.Range("E3").Value = 0.001558939920853764

$ws.Range("D4").Value = 0.1195411983538037
$ws.Range("E4").Value = 0.007524690390343425

$ws.Range("D5").Value = 0.140553776936452
$ws.Range("E5").Value = 0.001908700493080895

$ws.Range("D6").Value = 0.1361857591771688
$ws.Range("E6").Value = 0.007773459189339249

$ws.Range("D7").Value = 0.1460225114997454
$ws.Range("E7").Value = 0.01030732860520089

$ws.Range("D8").Value = 0.1273007915462319
$ws.Range("E8").Value = 0.008316008316008316

$ws.Range("D9").Value = 0.1276640149408169
$ws.Range("E9").Value = 0.005531897469494673

$ws.Range("E10").Value = 0.003037028238218831

# Restore sheet protection (no password is recoverable from the stored hash,
# so the sheet is re-protected with the same allow-list it shipped with).
$ws.Protect()
